$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N ("Late" moves N -> O,
# "Outstanding" moves O/P -> P/Q), matching the "Variable Instalments" column
# added to the Repayment Schedule sheet.
$ws.Columns("N:N").Insert() | Out-Null

# The newly inserted column gets a custom width (matches the author's resize).
$ws.Columns("N:N").ColumnWidth = 9.1666666666666666

# Make "Repayment Schedule" the active sheet (was "Transactions" before),
# and leave the selection on cell R8 as the author did.
$ws.Activate() | Out-Null
$ws.Range("R8").Select() | Out-Null
